# Continue the "Arbeitsjournal" table with a new entry for the next
# work session (started work on the next Carrier Board).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Add a new row at the bottom of the table; Word initializes it by
# cloning the formatting of the row immediately above it.
$newRow = $t.Rows.Add()
$rowIndex = $t.Rows.Count

$t.Cell($rowIndex, 1).Range.Text = "24.02.2025"
$t.Cell($rowIndex, 2).Range.Text = "120 min."
$t.Cell($rowIndex, 3).Range.Text = "D, P"
$t.Cell($rowIndex, 4).Range.Text = "Programmierung des ESP32, Entwurf eines PCBs, das alle nötigen Komponente halten kann"
